$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Marie Curie"
$ws.Range("B5").Value = "Otherland"
$ws.Range("C5").Value = "Rue de la Chanson"
$ws.Range("D5").Value = 43828
$ws.Range("E5").Value = "Paris"
$ws.Range("F5").Value = "Madame Curie"
